# Add a "traintype" column (I) that classifies each departure row as
# "spr" (Sprinter, stops at every listed station) or "IC" (Intercity,
# skips some stations).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "spr"
$ws.Range("I3").Value = "IC"
$ws.Range("I4").Value = "IC"
$ws.Range("I5").Value = "IC"
$ws.Range("I6").Value = "spr"
$ws.Range("I7").Value = "IC"
$ws.Range("I8").Value = "IC"
$ws.Range("I9").Value = "IC"
$ws.Range("I10").Value = "spr"
$ws.Range("I11").Value = "IC"
$ws.Range("I12").Value = "IC"
$ws.Range("I13").Value = "IC"
$ws.Range("I14").Value = "spr"
$ws.Range("I15").Value = "IC"
$ws.Range("I16").Value = "IC"
$ws.Range("I17").Value = "IC"

$ws.Range("I1").Value = "traintype"

$ws.Range("J4").Select() | Out-Null
